# Append 4 new data rows (8-11) to the IBBNamed worksheet, matching the
# "updated data" described in the commit. The sheet already has a header
# row (row 1) plus data rows 2-7; this adds rows 8-11 with the same shape,
# including the shared "Named" string in column N.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = 42613.76189814815
$ws.Range("B8").Value = 6
$ws.Range("C8").Value = 53
$ws.Range("D8").Value = 44
$ws.Range("E8").Value = 53
$ws.Range("F8").Value = 50
$ws.Range("G8").Value = 23615
$ws.Range("H8").Value = 11205
$ws.Range("I8").Value = 584
$ws.Range("J8").Value = 145
$ws.Range("K8").Value = 121
$ws.Range("L8").Value = 2
$ws.Range("M8").Value = 2
$ws.Range("N8").Value = "Named"

# Row 9
$ws.Range("A9").Value = 42613.891898148147
$ws.Range("B9").Value = 8
$ws.Range("C9").Value = 53
$ws.Range("D9").Value = 45
$ws.Range("E9").Value = 53
$ws.Range("F9").Value = 33
$ws.Range("G9").Value = 12986
$ws.Range("H9").Value = 10186
$ws.Range("I9").Value = 536
$ws.Range("J9").Value = 125
$ws.Range("K9").Value = 106
$ws.Range("L9").Value = 2
$ws.Range("M9").Value = 1
$ws.Range("N9").Value = "Named"

# Row 10
$ws.Range("A10").Value = 42614.88925925926
$ws.Range("B10").Value = 22
$ws.Range("C10").Value = 62
$ws.Range("D10").Value = 37
$ws.Range("E10").Value = 62
$ws.Range("F10").Value = 25
$ws.Range("G10").Value = 25883
$ws.Range("H10").Value = 19911
$ws.Range("I10").Value = 850
$ws.Range("J10").Value = 259
$ws.Range("K10").Value = 157
$ws.Range("L10").Value = 6
$ws.Range("M10").Value = 2
$ws.Range("N10").Value = "Named"

# Row 11
$ws.Range("A11").Value = 42615.887824074074
$ws.Range("B11").Value = 26
$ws.Range("C11").Value = 57
$ws.Range("D11").Value = 40
$ws.Range("E11").Value = 57
$ws.Range("F11").Value = 10
$ws.Range("G11").Value = 20601
$ws.Range("H11").Value = 13413
$ws.Range("I11").Value = 718
$ws.Range("J11").Value = 181
$ws.Range("K11").Value = 126
$ws.Range("L11").Value = 9
$ws.Range("M11").Value = 1
$ws.Range("N11").Value = "Named"
